$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  A = "Rio Grande do Sul"; B = "Diferença 2022-2000"; C = -0.136 },
    @{ Row = 3;  A = "Distrito Federal";  B = "Diferença 2022-2000"; C = -0.147 },
    @{ Row = 4;  A = "Rio de Janeiro";    B = "Diferença 2022-2000"; C = -0.151 },
    @{ Row = 5;  A = "Santa Catarina";    B = "Diferença 2022-2000"; C = -0.179 },
    @{ Row = 6;  A = "São Paulo";         B = "Diferença 2022-2000"; C = -0.179 },
    @{ Row = 7;  A = "Paraná";            B = "Diferença 2022-2000"; C = -0.201 },
    @{ Row = 8;  A = "Sergipe";           B = "Diferença 2022-2000"; C = -0.251 },
    @{ Row = 9;  A = "Nordeste";          B = "Diferença 2022-2000"; C = -0.248 },
    @{ Row = 10; A = "Brasil";            B = "Diferença 2022-2000"; C = -0.233 }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
}
